$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Grant the keys to new holders - append names to column B next to the
# existing numbered rows (13-17 correspond to key numbers 12-16).
# Order matches the original edit session so new shared-string entries land
# at the same table positions as the authored workbook.
$ws.Range("B13").Value = "郭泽民"
$ws.Range("B14").Value = "邱晨晨"
$ws.Range("B16").Value = "刘静"
$ws.Range("B17").Value = "赵燕"

# Update the summary cell (E1, merged E1:G1): used count 11 -> 16, spare 8 -> 3
$ws.Range("E1").Value = "（目前共19把，使用16把，闲置3把）"

$ws.Range("B15").Value = "张皓旭"

# Match the updated selection/active cell recorded in the saved view state
$ws.Range("E1:G1").Select()
